$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.804102
$ws.Range("H2").Value = 14.412306
$ws.Range("I2").Value = 0.2049869746002892
$ws.Range("J2").Value = 0.2049869746002892
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.238502
$ws.Range("N2").Value = 3.715506
$ws.Range("O2").Value = 0.2117788764206845
$ws.Range("P2").Value = 0.2117788764206845
$ws.Range("Q2").Value = 5.949889935204
$ws.Range("R2").Value = 53.549009416836
$ws.Range("S2").Value = 0.04341191116172465
$ws.Range("T2").Value = 0.04341191116172465

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.804102
$ws.Range("H3").Value = 14.412306
$ws.Range("I3").Value = 0.2049869746002892
$ws.Range("J3").Value = 0.2049869746002892
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.695367333333333
$ws.Range("N3").Value = 5.086101999999999
$ws.Range("O3").Value = 0.2899010166908616
$ws.Range("P3").Value = 0.2899010166908616
$ws.Range("Q3").Value = 8.144717596801334
$ws.Range("R3").Value = 73.302458371212
$ws.Range("S3").Value = 0.05942593234500767
$ws.Range("T3").Value = 0.05942593234500766

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.804102
$ws.Range("H4").Value = 14.412306
$ws.Range("I4").Value = 0.2049869746002892
$ws.Range("J4").Value = 0.2049869746002892
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.167891333333333
$ws.Range("N4").Value = 6.503674
$ws.Range("O4").Value = 0.3707007261800733
$ws.Range("P4").Value = 0.3707007261800732
$ws.Range("Q4").Value = 10.41477109024934
$ws.Range("R4").Value = 93.73293981224401
$ws.Range("S4").Value = 0.07598882034178345
$ws.Range("T4").Value = 0.07598882034178343

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.804102
$ws.Range("H5").Value = 14.412306
$ws.Range("I5").Value = 0.2049869746002892
$ws.Range("J5").Value = 0.2049869746002892
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7463296666666667
$ws.Range("N5").Value = 2.238989
$ws.Range("O5").Value = 0.1276193807083805
$ws.Range("P5").Value = 0.1276193807083805
$ws.Range("Q5").Value = 3.585443844292667
$ws.Range("R5").Value = 32.268994598634
$ws.Range("S5").Value = 0.02616031075177344
$ws.Range("T5").Value = 0.02616031075177343

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.489274999999999
$ws.Range("H6").Value = 25.467825
$ws.Range("I6").Value = 0.3622301938634671
$ws.Range("J6").Value = 0.3622301938634671
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.238502
$ws.Range("N6").Value = 3.715506
$ws.Range("O6").Value = 0.2117788764206845
$ws.Range("P6").Value = 0.2117788764206845
$ws.Range("Q6").Value = 10.51398406605
$ws.Range("R6").Value = 94.62585659444998
$ws.Range("S6").Value = 0.07671270346205181
$ws.Range("T6").Value = 0.07671270346205179

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.489274999999999
$ws.Range("H7").Value = 25.467825
$ws.Range("I7").Value = 0.3622301938634671
$ws.Range("J7").Value = 0.3622301938634671
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.695367333333333
$ws.Range("N7").Value = 5.086101999999999
$ws.Range("O7").Value = 0.2899010166908616
$ws.Range("P7").Value = 0.2899010166908616
$ws.Range("Q7").Value = 14.39243951868333
$ws.Range("R7").Value = 129.53195566815
$ws.Range("S7").Value = 0.105010901477147
$ws.Range("T7").Value = 0.105010901477147

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.489274999999999
$ws.Range("H8").Value = 25.467825
$ws.Range("I8").Value = 0.3622301938634671
$ws.Range("J8").Value = 0.3622301938634671
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.167891333333333
$ws.Range("N8").Value = 6.503674
$ws.Range("O8").Value = 0.3707007261800733
$ws.Range("P8").Value = 0.3707007261800732
$ws.Range("Q8").Value = 18.40382569878333
$ws.Range("R8").Value = 165.63443128905
$ws.Range("S8").Value = 0.134278995909536
$ws.Range("T8").Value = 0.1342789959095359

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.489274999999999
$ws.Range("H9").Value = 25.467825
$ws.Range("I9").Value = 0.3622301938634671
$ws.Range("J9").Value = 0.3622301938634671
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7463296666666667
$ws.Range("N9").Value = 2.238989
$ws.Range("O9").Value = 0.1276193807083805
$ws.Range("P9").Value = 0.1276193807083805
$ws.Range("Q9").Value = 6.335797780991666
$ws.Range("R9").Value = 57.022180028925
$ws.Range("S9").Value = 0.04622759301473229
$ws.Range("T9").Value = 0.04622759301473228

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.319169666666667
$ws.Range("H10").Value = 6.957509
$ws.Range("I10").Value = 0.09895701081175237
$ws.Range("J10").Value = 0.09895701081175236
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.238502
$ws.Range("N10").Value = 3.715506
$ws.Range("O10").Value = 0.2117788764206845
$ws.Range("P10").Value = 0.2117788764206845
$ws.Range("Q10").Value = 2.872296270506
$ws.Range("R10").Value = 25.850666434554
$ws.Range("S10").Value = 0.02095700456366245
$ws.Range("T10").Value = 0.02095700456366245

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.319169666666667
$ws.Range("H11").Value = 6.957509
$ws.Range("I11").Value = 0.09895701081175237
$ws.Range("J11").Value = 0.09895701081175236
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.695367333333333
$ws.Range("N11").Value = 5.086101999999999
$ws.Range("O11").Value = 0.2899010166908616
$ws.Range("P11").Value = 0.2899010166908616
$ws.Range("Q11").Value = 3.931844493324222
$ws.Range("R11").Value = 35.386600439918
$ws.Range("S11").Value = 0.0286877380430156
$ws.Range("T11").Value = 0.02868773804301559

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.319169666666667
$ws.Range("H12").Value = 6.957509
$ws.Range("I12").Value = 0.09895701081175237
$ws.Range("J12").Value = 0.09895701081175236
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.167891333333333
$ws.Range("N12").Value = 6.503674
$ws.Range("O12").Value = 0.3707007261800733
$ws.Range("P12").Value = 0.3707007261800732
$ws.Range("Q12").Value = 5.027707820896222
$ws.Range("R12").Value = 45.249370388066
$ws.Range("S12").Value = 0.03668343576852597
$ws.Range("T12").Value = 0.03668343576852596

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.319169666666667
$ws.Range("H13").Value = 6.957509
$ws.Range("I13").Value = 0.09895701081175237
$ws.Range("J13").Value = 0.09895701081175236
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.7463296666666667
$ws.Range("N13").Value = 2.238989
$ws.Range("O13").Value = 0.1276193807083805
$ws.Range("P13").Value = 0.1276193807083805
$ws.Range("Q13").Value = 1.730865124266778
$ws.Range("R13").Value = 15.577786118401
$ws.Range("S13").Value = 0.01262883243654835
$ws.Range("T13").Value = 0.01262883243654835

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 7.823586333333334
$ws.Range("H14").Value = 23.470759
$ws.Range("I14").Value = 0.3338258207244912
$ws.Range("J14").Value = 0.3338258207244912
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.238502
$ws.Range("N14").Value = 3.715506
$ws.Range("O14").Value = 0.2117788764206845
$ws.Range("P14").Value = 0.2117788764206845
$ws.Range("Q14").Value = 9.689527321006
$ws.Range("R14").Value = 87.205745889054
$ws.Range("S14").Value = 0.07069725723324563
$ws.Range("T14").Value = 0.07069725723324563

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 7.823586333333334
$ws.Range("H15").Value = 23.470759
$ws.Range("I15").Value = 0.3338258207244912
$ws.Range("J15").Value = 0.3338258207244912
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.695367333333333
$ws.Range("N15").Value = 5.086101999999999
$ws.Range("O15").Value = 0.2899010166908616
$ws.Range("P15").Value = 0.2899010166908616
$ws.Range("Q15").Value = 13.26385269904644
$ws.Range("R15").Value = 119.374674291418
$ws.Range("S15").Value = 0.09677644482569131
$ws.Range("T15").Value = 0.09677644482569131

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 7.823586333333334
$ws.Range("H16").Value = 23.470759
$ws.Range("I16").Value = 0.3338258207244912
$ws.Range("J16").Value = 0.3338258207244912
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2.167891333333333
$ws.Range("N16").Value = 6.503674
$ws.Range("O16").Value = 0.3707007261800733
$ws.Range("P16").Value = 0.3707007261800732
$ws.Range("Q16").Value = 16.96068500761845
$ws.Range("R16").Value = 152.646165068566
$ws.Range("S16").Value = 0.1237494741602279
$ws.Range("T16").Value = 0.1237494741602278

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 7.823586333333334
$ws.Range("H17").Value = 23.470759
$ws.Range("I17").Value = 0.3338258207244912
$ws.Range("J17").Value = 0.3338258207244912
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.7463296666666667
$ws.Range("N17").Value = 2.238989
$ws.Range("O17").Value = 0.1276193807083805
$ws.Range("P17").Value = 0.1276193807083805
$ws.Range("Q17").Value = 5.838974580294556
$ws.Range("R17").Value = 52.55077122265101
$ws.Range("S17").Value = 0.04260264450532644
$ws.Range("T17").Value = 0.04260264450532644
